$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D historically stores prices as text (inline strings), even when
# the values look numeric (e.g. "51.048.36" uses "." as a thousands separator).
# Force text format on column D cells before assigning so Excel does not
# silently convert them into numeric values.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.048.36'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.957.35'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '380.28'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.25'
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('E7').Value = '  +1.06%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.588'
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.35'
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0851'
$ws.Range('E12').Value = '  +1.75%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.48'
$ws.Range('E13').Value = '  +2.64%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.423.90'
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.77'
$ws.Range('E15').Value = '  +5.74%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '12.38'
$ws.Range('E16').Value = '  +74.75%  '
$ws.Range('E17').Value = '  +4.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.961.26'
$ws.Range('E18').Value = '  +1.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '51.056.12'
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('E20').Value = '  -1.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.40'
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0955'
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.36'
$ws.Range('E23').Value = '  +16.88%  '
$ws.Range('E24').Value = '  +2.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '266.66'
$ws.Range('E25').Value = '  +1.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.06'
$ws.Range('E26').Value = '  -2.01%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.86'
$ws.Range('E28').Value = '  +1.31%  '
$ws.Range('E29').Value = '  -1.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.98'
$ws.Range('E30').Value = '  -8.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.107'
$ws.Range('E31').Value = '  -6.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.38'
$ws.Range('E32').Value = '  +6.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.11'
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.02'
$ws.Range('E35').Value = '  -1.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0435'
$ws.Range('E36').Value = '  -3.51%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.19'
$ws.Range('E38').Value = '  +7.74%  '
$ws.Range('E39').Value = '  +1.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.65'
$ws.Range('E40').Value = '  +0.74%  '
$ws.Range('E41').Value = '  +3.25%  '
$ws.Range('E42').Value = '  -3.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '118.42'
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.57'
$ws.Range('E44').Value = '  +11.49%  '
$ws.Range('E45').Value = '  +1.53%  '
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.024.42'
$ws.Range('E47').Value = '  +1.21%  '
$ws.Range('E48').Value = '  -2.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.256'
$ws.Range('E49').Value = '  -5.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0318'
$ws.Range('E50').Value = '  -7.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.35'
$ws.Range('E51').Value = '  +6.67%  '
